$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$aVals = @(0,1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29)
$bVals = @(40,40,40,40,40,40,40,45,51,55,60,63,67,71,74,77,80,83,86,89,92,95,98,101,104,108,111,115,119,124)
for ($i = 0; $i -lt $aVals.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $aVals[$i]
    $ws.Cells.Item($row, 2).Value = $bVals[$i]
}

$ws = $wb.Worksheets.Item(2)
$aVals = @(0,1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29)
$bVals = @(40,40,40,40,40,40,40,45,51,56,60,64,67,71,74,77,80,83,86,89,92,95,98,101,104,107,111,114,118,123)
for ($i = 0; $i -lt $aVals.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $aVals[$i]
    $ws.Cells.Item($row, 2).Value = $bVals[$i]
}

$ws = $wb.Worksheets.Item(3)
$aVals = @(0,1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29)
$bVals = @(40,40,40,40,40,40,40,46,52,56,61,64,68,71,75,78,81,84,87,89,92,95,98,101,104,107,110,114,117,122)
for ($i = 0; $i -lt $aVals.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $aVals[$i]
    $ws.Cells.Item($row, 2).Value = $bVals[$i]
}

$ws = $wb.Worksheets.Item(4)
$aVals = @(0,1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29)
$bVals = @(40,40,40,40,40,40,42,48,53,58,62,66,69,72,76,79,82,85,87,90,93,96,99,101,104,107,110,114,117,121)
for ($i = 0; $i -lt $aVals.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $aVals[$i]
    $ws.Cells.Item($row, 2).Value = $bVals[$i]
}

$ws = $wb.Worksheets.Item(5)
$aVals = @(0,1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29)
$bVals = @(40,40,40,40,40,40,46,52,56,61,64,68,71,75,78,81,84,86,89,92,95,97,100,103,106,108,111,114,118,121)
for ($i = 0; $i -lt $aVals.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $aVals[$i]
    $ws.Cells.Item($row, 2).Value = $bVals[$i]
}

$ws = $wb.Worksheets.Item(6)
$aVals = @(0,1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29)
$bVals = @(40,40,40,40,42,50,55,60,64,68,71,75,78,81,83,86,89,92,94,97,99,102,105,107,110,112,115,118,121,124)
for ($i = 0; $i -lt $aVals.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $aVals[$i]
    $ws.Cells.Item($row, 2).Value = $bVals[$i]
}
